$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the second data record (row 3) entirely
$ws.Range("A3").EntireRow.Delete()

# Remove "Sudah Pengumuman" (E) and "Tanggal Pendataan" (F) columns
$ws.Range("E:F").Delete()

# Update remaining data row with corrected values
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = '001'
$ws.Range("C2").Value = '0011'
$ws.Range("E2").Value = 12345
$ws.Range("F2").Value = 12345
$ws.Range("G2").Value = 11.22
$ws.Range("H2").Value = '012345'
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = '10%'
$ws.Range("L2").Value = '3509072008000002'
$ws.Range("M2").Value = 'SUDRAJAD HADI SAPUTRA'
$ws.Range("N2").Value = 'JEMBER'
$ws.Range("O2").Value = '2000-08-20'
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 'JALAN BENDO GANG 3 NO. 36 DESA SIDOMEKAR KECAMATAN SEMBORO'
$ws.Range("R2").Value = 'ISLAM'
$ws.Range("S2").Value = 'MAHASISWA'
$ws.Range("T2").Value = '3509072008000002'
$ws.Range("U2").Value = 'SUDRAJAD HADI SAPUTRA'
$ws.Range("V2").Value = 'JEMBER'
$ws.Range("W2").Value = '2000-08-20'
$ws.Range("X2").Value = 22
$ws.Range("Y2").Value = 'JALAN BENDO GANG 3 NO. 36 DESA SIDOMEKAR KECAMATAN SEMBORO'
$ws.Range("Z2").Value = '002'
$ws.Range("AA2").Value = '003'
$ws.Range("AB2").Value = 'SONGON'
$ws.Range("AC2").Value = 'PONDOKJOYO'
$ws.Range("AD2").Value = 'SEMBORO'
$ws.Range("AE2").Value = 2000
$ws.Range("AF2").Value = 12345
$ws.Range("AG2").Value = 12345
$ws.Range("AH2").Value = 'S. II'
$ws.Range("AI2").Value = 'YASAN'
$ws.Range("AJ2").Value = 'PEKARANGAN'
$ws.Range("AK2").Value = 900
$ws.Range("AL2").Value = 'JALAN DESA'
$ws.Range("AM2").Value = 'TANAH ORANG'
$ws.Range("AN2").Value = 'TANAH ORANG'
$ws.Range("AO2").Value = 'JALAN DESA'
$ws.Range("AP2").Value = 2000
$ws.Range("AQ2").Value = 'ORANG 1'
$ws.Range("AR2").Value = 2010
$ws.Range("AS2").Value = 'ORANG 2'
$ws.Range("AT2").Value = 'WARIS'
$ws.Range("AU2").Value = 'WARIS 2'
$ws.Range("AV2").Value = 'ORANG 2'
$ws.Range("AW2").Value = 2023
$ws.Range("AX2").Value = ""
$ws.Range("AY2").Value = 'SUDRAJAD HADI SAPUTRA'
$ws.Range("BB2").Value = ""
$ws.Range("BC2").Value = 'JUAL BELI 3'
$ws.Range("BE2").Value = 'JUAL BELI 3'
$ws.Range("BF2").Value = 'DIDIK SAENULLA'
$ws.Range("BS2").Value = 'EDI BEDOR'
